$d = $word.ActiveDocument

# Change 1: update the Celular/Telegram phone line text
$d.Content.Find.Execute(
    "Celular/Telegram: (77) 98139-8699",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Celular: (77) 98139-8699 (WhatsApp e Telegram)",
    2
) | Out-Null
